$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.7791253856517797, 1.1073232676938676, 1.2645139860528265, 2.7786174103803853)"
$ws.Range("C2").Value = "NIG(1.3366693574265214, 1.0508698382454538, 3.7029111253614975, 5.653136882029561)"
$ws.Range("D2").Value = "JSU(-1.0345121350878124, 1.0579872079252002, -0.25716504021565123, 2.6109531569470708)"
$ws.Range("E2").Value = "NIG(1.3873503671575445, 1.0968709879154503, 3.4152504510461563, 5.242283163666832)"
